$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + volume(1h)) and the Quant/FraxShare row swap
# Row 2
$ws.Range("D2").Value = '26.880.24'
$ws.Range("E2").Value = '  -0.41%  '
# Row 3
$ws.Range("D3").Value = '1.859.42'
$ws.Range("E3").Value = '  -0.08%  '
# Row 4
$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  -0.14%  '
# Row 5
$ws.Range("D5").Value = '304.62'
# Row 6
$ws.Range("D6").Value = '''1.000'
$ws.Range("E6").Value = '  -0.12%  '
# Row 7
$ws.Range("E7").Value = '  -1.49%  '
# Row 8
$ws.Range("E8").Value = '  -2.38%  '
# Row 9
$ws.Range("E9").Value = '  +0.99%  '
# Row 10
$ws.Range("D10").Value = '''0.8920'
$ws.Range("E10").Value = '  +0.50%  '
# Row 11
$ws.Range("D11").Value = '20.68'
$ws.Range("E11").Value = '  +0.59%  '
# Row 12
$ws.Range("D12").Value = '1.875.58'
$ws.Range("E12").Value = '  +0.82%  '
# Row 13
$ws.Range("D13").Value = '0.07475'
$ws.Range("E13").Value = '  -0.95%  '
# Row 14
$ws.Range("D14").Value = '93.82'
$ws.Range("E14").Value = '  +5.92%  '
# Row 15
$ws.Range("D15").Value = '5.226'
# Row 16
$ws.Range("E16").Value = '  -0.11%  '
# Row 17
$ws.Range("D17").Value = '''0.000008494'
$ws.Range("E17").Value = '  +1.05%  '
# Row 18
$ws.Range("E18").Value = '  +0.81%  '
# Row 19
$ws.Range("D19").Value = '''1.000'
$ws.Range("E19").Value = '  -0.20%  '
# Row 20
$ws.Range("D20").Value = '26.933.83'
$ws.Range("E20").Value = '  -0.40%  '
# Row 21
$ws.Range("D21").Value = '5.023'
$ws.Range("E21").Value = '  -0.59%  '
# Row 22
$ws.Range("D22").Value = '2.115.51'
$ws.Range("E22").Value = '  +0.85%  '
# Row 23
$ws.Range("D23").Value = '10.38'
$ws.Range("E23").Value = '  -1.41%  '
# Row 24
$ws.Range("D24").Value = '6.416'
$ws.Range("E24").Value = '  -0.90%  '
# Row 25
$ws.Range("D25").Value = '147.62'
$ws.Range("E25").Value = '  -1.15%  '
# Row 26
$ws.Range("D26").Value = '1.785'
$ws.Range("E26").Value = '  -2.97%  '
# Row 27
$ws.Range("D27").Value = '17.85'
$ws.Range("E27").Value = '  -0.61%  '
# Row 28
$ws.Range("D28").Value = '2.082'
$ws.Range("E28").Value = '  -0.88%  '
# Row 29
$ws.Range("D29").Value = '113.01'
$ws.Range("E29").Value = '  +0.22%  '
# Row 30
$ws.Range("D30").Value = '4.687'
$ws.Range("E30").Value = '  +0.15%  '
# Row 31
$ws.Range("D31").Value = '''4.670'
$ws.Range("E31").Value = '  +0.44%  '
# Row 32
$ws.Range("D32").Value = '0.09214'
$ws.Range("E32").Value = '  +2.02%  '
# Row 33
$ws.Range("D33").Value = '0.05122'
$ws.Range("E33").Value = '  +0.12%  '
# Row 34
$ws.Range("D34").Value = '0.7472'
$ws.Range("E34").Value = '  +2.22%  '
# Row 35
$ws.Range("D35").Value = '2.972'
$ws.Range("E35").Value = '  -3.33%  '
# Row 36
$ws.Range("D36").Value = '1.151'
$ws.Range("E36").Value = '  -0.12%  '
# Row 37
$ws.Range("D37").Value = '3.258'
$ws.Range("E37").Value = '  +6.97%  '
# Row 38
$ws.Range("D38").Value = '2.572'
$ws.Range("E38").Value = '  +3.21%  '
# Row 39
$ws.Range("E39").Value = '  -2.28%  '
# Row 40
$ws.Range("D40").Value = '0.5556'
$ws.Range("E40").Value = '  +4.57%  '
# Row 41
$ws.Range("D41").Value = '1.071'
$ws.Range("E41").Value = '  -0.16%  '
# Row 42
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '6.533'
$ws.Range("E42").Value = '  -0.82%  '
# Row 43
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '117.13'
$ws.Range("E43").Value = '  +1.40%  '
# Row 44
$ws.Range("D44").Value = '8.503'
$ws.Range("E44").Value = '  +2.66%  '
# Row 45
$ws.Range("D45").Value = '0.1469'
$ws.Range("E45").Value = '  -0.03%  '
# Row 46
$ws.Range("D46").Value = '0.4673'
$ws.Range("E46").Value = '  +1.20%  '
# Row 47
$ws.Range("D47").Value = '''1.000'
$ws.Range("E47").Value = '  -0.18%  '
# Row 48
$ws.Range("D48").Value = '10.04'
$ws.Range("E48").Value = '  -0.05%  '
# Row 49
$ws.Range("D49").Value = '1.559'
$ws.Range("E49").Value = '  -0.26%  '
# Row 50
$ws.Range("D50").Value = '''36.70'
$ws.Range("E50").Value = '  +0.10%  '
# Row 51
$ws.Range("D51").Value = '62.94'
$ws.Range("E51").Value = '  -1.85%  '
